# Trade #11 closed at 2026-02-16 22:58:23 - base_strategy DOWN +0.000%
# Appends a new trade row (row 12) to both the "All Trades" and
# "base_strategy" sheets, mirroring the existing OPEN-trade row layout.

$wb = $excel.ActiveWorkbook

$row = 12

$tradeNum    = 11
$tradeDate   = "2026-02-16"
$tradeTime   = "22:58:23"
$strategy    = "base_strategy"
$side        = "DOWN"
$entryPrice  = 0.5
$status      = "OPEN"
$pnlPct      = 0
$pnlUsd      = 0
$capAfter    = 100
$entrySlip   = 0
$exitSlip    = 0
$confidence  = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$duration    = 0

foreach ($sheetName in @("All Trades", "base_strategy")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Trade #
    $ws.Cells.Item($row, 1).Value = $tradeNum

    # Date - force text storage so Excel doesn't reinterpret "2026-02-16"
    # as a date serial number (write with quote-prefix, then clear the
    # quote-prefix style so the stored cell carries no formatting).
    $ws.Cells.Item($row, 2).Value = "'" + $tradeDate
    $ws.Cells.Item($row, 2).Style = "Normal"

    # Time (plain text, Excel keeps this as text already)
    $ws.Cells.Item($row, 3).Value = $tradeTime

    # Strategy / Side
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side

    # Entry Price
    $ws.Cells.Item($row, 6).Value = $entryPrice

    # Exit Price - empty text cell (trade still OPEN). A bare "" assignment
    # collapses to a true blank cell, so enter a lone quote-prefix (empty
    # text) and then strip the quote-prefix style, leaving an empty text
    # value with no formatting - matching the other OPEN rows above it.
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    # Status
    $ws.Cells.Item($row, 8).Value = $status

    # P&L %, P&L $, Capital After
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlUsd
    $ws.Cells.Item($row, 11).Value = $capAfter

    # Entry / Exit slippage
    $ws.Cells.Item($row, 12).Value = $entrySlip
    $ws.Cells.Item($row, 13).Value = $exitSlip

    # Confidence
    $ws.Cells.Item($row, 14).Value = $confidence

    # Entry Reason
    $ws.Cells.Item($row, 15).Value = $entryReason

    # Exit Reason - empty text cell (same trick as Exit Price above)
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).Style = "Normal"

    # Duration (min)
    $ws.Cells.Item($row, 17).Value = $duration
}
